$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.558.41"
$ws.Range("E2").Value = "  +2.53%  "
$ws.Range("D3").Value = "1.671.01"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "239.09"
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "0.4793"
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("E8").Value = "  +3.20%  "
$ws.Range("E9").Value = "  +3.20%  "
$ws.Range("E10").Value = "  -2.51%  "
$ws.Range("D11").Value = "1.670.39"
$ws.Range("E11").Value = "  +2.10%  "
$ws.Range("D12").Value = "14.88"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "0.5916"
$ws.Range("E13").Value = "  -3.69%  "
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").Value = "75.37"
$ws.Range("E15").Value = "  +3.98%  "
$ws.Range("D16").Value = "0.9999"
$ws.Range("D17").Value = "0.9997"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "25.548.56"
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("D19").Value = "0.000006773"
$ws.Range("E19").Value = "  +2.91%  "
$ws.Range("D20").Value = "11.49"
$ws.Range("E20").Value = "  +2.03%  "
$ws.Range("D21").Value = "1.883.73"
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "8.749"
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("D24").Value = "5.286"
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").Value = "136.77"
$ws.Range("E25").Value = "  +3.57%  "
$ws.Range("D26").Value = "15.08"
$ws.Range("E26").Value = "  +2.05%  "
$ws.Range("D27").Value = "1.392"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").Value = "1.727"
$ws.Range("E28").Value = "  +4.50%  "
$ws.Range("D29").Value = "104.87"
$ws.Range("D30").Value = "3.977"
$ws.Range("E30").Value = "  +6.94%  "
$ws.Range("D31").Value = "0.07838"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").Value = "3.657"
$ws.Range("E32").Value = "  +3.63%  "
$ws.Range("D33").Value = "0.9989"
$ws.Range("E34").Value = "  -3.95%  "
$ws.Range("D35").Value = "2.617"
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("D36").Value = "0.6110"
$ws.Range("E36").Value = "  +5.05%  "
$ws.Range("D37").Value = "0.9518"
$ws.Range("E37").Value = "  +3.22%  "
$ws.Range("D38").Value = "2.596"
$ws.Range("E38").Value = "  +2.32%  "
$ws.Range("D39").Value = "0.8567"
$ws.Range("E39").Value = "  +2.40%  "
$ws.Range("D40").Value = "0.9996"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").Value = "1.875"
$ws.Range("E41").Value = "  +4.00%  "
$ws.Range("D42").Value = "0.01475"
$ws.Range("E42").Value = "  -5.15%  "
$ws.Range("D43").Value = "96.21"
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("D44").Value = "0.3780"
$ws.Range("E44").Value = "  +2.29%  "
$ws.Range("D45").Value = "4.865"
$ws.Range("E45").Value = "  +2.47%  "
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("D47").Value = "6.230"
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("D48").Value = "0.05262"
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("D49").Value = "29.86"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("D50").Value = "7.392"
$ws.Range("E50").Value = "  +2.83%  "
$ws.Range("D51").Value = "1.001"
$ws.Range("E51").Value = "  +0.19%  "
